# "Update code, fix logic"
#
# Sheet "1" (table Tabela13 / Tabela135 shape): fill in the two
# "Número da conta" cells that were previously blank, and correct the
# procedure-code typo in B2.
#
# Sheet "2" (table Tabela1343): correct the procedure code and the two
# "Valor unitário" values.
#
# The writes are issued in this specific order so the shared-string
# table grows in the same sequence as the authoritative edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1")
$ws2 = $wb.Worksheets.Item("2")

$ws1.Range("A2").Value = "22222222"
$ws1.Range("A3").Value = "22106869"
$ws1.Range("E3").Value = "10"

$ws2.Range("C2").Value = "11,53"
$ws2.Range("D2").Value = "12,00"
$ws2.Range("B2").Value = "asfsaf"

$ws1.Range("B2").Value = "111"

# Sheet "2" was the active/selected tab in the original file (selection
# left sitting on C3); the edit session ends with sheet "1" focused and
# both sheets' cursors resting on B2.
[void]$ws2.Range("B2").Select()
[void]$ws1.Activate()
[void]$ws1.Range("B2").Select()
